$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the two existing rows (A2:B3) down into the
# three new rows (A4:B6), reusing the existing cell styles exactly.
$ws.Range("A2:B3").Copy()
$ws.Range("A4:B6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A4").Value = 45183
$ws.Range("B4").Value = 0.4069444444444445

$ws.Range("A5").Value = 45184
$ws.Range("B5").Value = 0.40972222222222227

$ws.Range("A6").Value = 45185
$ws.Range("B6").Value = 0.41666666666666669

[void]$ws.Range("A7").Select()
